# Regenerate save_data to use K instead of Strike# (column G, header "K")
# Update the K values for each row based on recalculated stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 1
    4  = 2
    5  = 1
    6  = 2
    7  = 0
    8  = 1
    9  = 2
    10 = 1
    11 = 3
    12 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
